# Generate Report for Handoff
# Update the localization-status report: file "b.md" has now been handed off
# (new handoff xliff files generated) instead of previously being
# "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"
$newHandoffDate = "2016-08-29 10:37:14"

# ----------------------------------------------------------------------
# "Overview" sheet: row 3 is the b.md file.
# ----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusReady
$wsOverview.Range("F3").Value = $statusReady
$wsOverview.Range("G3").Value = $newHandoffDate

# ----------------------------------------------------------------------
# "zh-cn" sheet: row 3 is the b.md file.
# ----------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusReady
# Prefix with an apostrophe so the literal text "False" is not auto-coerced
# into a Boolean cell, then clear the resulting quote-prefix formatting so
# the cell keeps its original (default) style.
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").ClearFormats()
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-29 10:37:09"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a407918c706a037fa3086e20325b1914914102a8/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1967caa7bc98a8e63a428ee0f88c610f6ebc566c/e2e/b.md."
$wsZhCn.Columns(16).ColumnWidth = 39.15

# ----------------------------------------------------------------------
# "de-de" sheet: row 3 is the b.md file.
# ----------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusReady
# Prefix with an apostrophe so the literal text "False" is not auto-coerced
# into a Boolean cell, then clear the resulting quote-prefix formatting so
# the cell keeps its original (default) style.
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").ClearFormats()
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = $newHandoffDate
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a407918c706a037fa3086e20325b1914914102a8/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1967caa7bc98a8e63a428ee0f88c610f6ebc566c/e2e/b.md."
$wsDeDe.Columns(16).ColumnWidth = 39.15
